# Auto-generated edit script
# Applies per-cell numeric updates (and a few cell clears) as described by the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 122.166664
$ws.Range("I39").Value = 140.6
$ws.Range("J39").Value = 30
$ws.Range("K39").Value = 421.8
$ws.Range("L39").Value = 90
$ws.Range("M39").Value = -125.8
$ws.Range("N39").Value = -682
$ws.Range("H86").Value = 4475.593
$ws.Range("I86").Value = 760.6429000000001
$ws.Range("J86").Value = 8476.308000000001
$ws.Range("K86").Value = 760.6429000000001
$ws.Range("L86").Value = 8476.308000000001
$ws.Range("M86").Value = 362.3570999999999
$ws.Range("N86").Value = -10722.308
$ws.Range("H89").Value = 4475.593
$ws.Range("I89").Value = 760.6429000000001
$ws.Range("J89").Value = 8476.308000000001
$ws.Range("K89").Value = 3803.2145
$ws.Range("L89").Value = 42381.54000000001
$ws.Range("M89").Value = 1812.7855
$ws.Range("N89").Value = -53613.54000000001
$ws.Range("H111").Value = 2939.5
$ws.Range("I111").Value = 3422.3076
$ws.Range("J111").Value = 2042.8572
$ws.Range("K111").Value = 10266.9228
$ws.Range("L111").Value = 6128.571599999999
$ws.Range("M111").Value = -7199.9228
$ws.Range("N111").Value = -12262.5716
$ws.Range("H132").Value = 4278.05
$ws.Range("I132").Value = 4278.05
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12834.15
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10304.15
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6639
$ws.Range("I32").Value = 6336.925
$ws.Range("K32").Value = 6336.925
$ws.Range("M32").Value = -6049.925
$ws.Range("H130").Value = 29776.924
$ws.Range("J130").Value = 29776.924
$ws.Range("L130").Value = 29776.924
$ws.Range("N130").Value = -39816.924
$ws.Range("H132").Value = 16238.828
$ws.Range("I132").Value = 1619.5518
$ws.Range("K132").Value = 4858.6554
$ws.Range("M132").Value = -2328.6554
$ws.Range("H138").Value = 48524.168
$ws.Range("J138").Value = 48524.168
$ws.Range("L138").Value = 48524.168
$ws.Range("N138").Value = -58804.168

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1498.5686
$ws.Range("J86").Value = 1676.1904
$ws.Range("L86").Value = 1676.1904
$ws.Range("N86").Value = -3922.1904
$ws.Range("H89").Value = 1498.5686
$ws.Range("J89").Value = 1676.1904
$ws.Range("L89").Value = 8380.951999999999
$ws.Range("N89").Value = -19612.952
$ws.Range("H134").Value = 5583.684
$ws.Range("I134").Value = 5583.684
$ws.Range("K134").Value = 16751.052
$ws.Range("M134").Value = -14216.052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4521.727
$ws.Range("J62").Value = 4251.5
$ws.Range("L62").Value = 4251.5
$ws.Range("N62").Value = -5499.5
$ws.Range("H65").Value = 4521.727
$ws.Range("J65").Value = 4251.5
$ws.Range("L65").Value = 21257.5
$ws.Range("N65").Value = -27497.5
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H132").Value = 4331.2
$ws.Range("I132").Value = 1700
$ws.Range("J132").Value = 6085.3335
$ws.Range("K132").Value = 5100
$ws.Range("L132").Value = 18256.0005
$ws.Range("M132").Value = -2570
$ws.Range("N132").Value = -23316.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 420.55554
$ws.Range("I23").Value = 50
$ws.Range("J23").Value = 526.4286
$ws.Range("K23").Value = 150
$ws.Range("L23").Value = 1579.2858
$ws.Range("M23").Value = 85
$ws.Range("N23").Value = -2049.2858
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H97").Value = 693.3333
$ws.Range("J97").Value = 693.3333
$ws.Range("L97").Value = 2079.9999
$ws.Range("N97").Value = -3071.9999
$ws.Range("H113").Value = 376.75
$ws.Range("J113").Value = 373.4375
$ws.Range("L113").Value = 1120.3125
$ws.Range("N113").Value = -5460.3125
$ws.Range("H131").Value = 714.62
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 725.9158
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 2177.7474
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -12257.7474
$ws.Range("H140").Value = 2508.5715
$ws.Range("J140").Value = 3827.1428
$ws.Range("L140").Value = 11481.4284
$ws.Range("N140").Value = -21841.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3481287.5
$ws.Range("I70").Value = 3877.7
$ws.Range("J70").Value = 7828050
$ws.Range("K70").Value = 3877.7
$ws.Range("L70").Value = 7828050
$ws.Range("M70").Value = -3607.7
$ws.Range("N70").Value = -7828590
$ws.Range("H73").Value = 3481287.5
$ws.Range("I73").Value = 3877.7
$ws.Range("J73").Value = 7828050
$ws.Range("K73").Value = 3877.7
$ws.Range("L73").Value = 7828050
$ws.Range("M73").Value = -2941.7
$ws.Range("N73").Value = -7829922
$ws.Range("H80").Value = 4255
$ws.Range("I80").Value = 3633.3333
$ws.Range("J80").Value = 4521.4287
$ws.Range("K80").Value = 3633.3333
$ws.Range("L80").Value = 4521.4287
$ws.Range("M80").Value = -2635.3333
$ws.Range("N80").Value = -6517.4287
$ws.Range("H83").Value = 4255
$ws.Range("I83").Value = 3633.3333
$ws.Range("J83").Value = 4521.4287
$ws.Range("K83").Value = 18166.6665
$ws.Range("L83").Value = 22607.1435
$ws.Range("M83").Value = -13174.6665
$ws.Range("N83").Value = -32591.1435
$ws.Range("H102").Value = 2821.913
$ws.Range("J102").Value = 2346.8572
$ws.Range("L102").Value = 2346.8572
$ws.Range("N102").Value = -5590.8572
$ws.Range("H113").Value = 2487.889
$ws.Range("I113").Value = 1979
$ws.Range("J113").Value = 3287.5715
$ws.Range("K113").Value = 1979
$ws.Range("L113").Value = 3287.5715
$ws.Range("M113").Value = 191
$ws.Range("N113").Value = -7627.5715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7000
$ws.Range("I22").Value = 5500.5
$ws.Range("J22").Value = 9999
$ws.Range("K22").Value = 5500.5
$ws.Range("L22").Value = 9999
$ws.Range("M22").Value = -5205.5
$ws.Range("N22").Value = -10589
$ws.Range("H27").Value = 7000
$ws.Range("I27").Value = 5500.5
$ws.Range("J27").Value = 9999
$ws.Range("K27").Value = 5500.5
$ws.Range("L27").Value = 9999
$ws.Range("M27").Value = -5393.5
$ws.Range("N27").Value = -10213
$ws.Range("H40").Value = 2958.9707
$ws.Range("I40").Value = 2332.5
$ws.Range("K40").Value = 2332.5
$ws.Range("M40").Value = -2196.5
$ws.Range("H61").Value = 4833.857
$ws.Range("I61").Value = 1910.7142
$ws.Range("J61").Value = 7757
$ws.Range("K61").Value = 1910.7142
$ws.Range("L61").Value = 7757
$ws.Range("M61").Value = -1708.7142
$ws.Range("N61").Value = -8161
$ws.Range("H105").Value = 23699.75
$ws.Range("J105").Value = 23699.75
$ws.Range("L105").Value = 23699.75
$ws.Range("N105").Value = -30687.75
$ws.Range("H113").Value = 4833.857
$ws.Range("I113").Value = 1910.7142
$ws.Range("J113").Value = 7757
$ws.Range("K113").Value = 1910.7142
$ws.Range("L113").Value = 7757
$ws.Range("M113").Value = 259.2858000000001
$ws.Range("N113").Value = -12097
$ws.Range("H132").Value = 3684.4443
$ws.Range("I132").Value = 3218.6667
$ws.Range("J132").Value = 4616
$ws.Range("K132").Value = 9656.000100000001
$ws.Range("L132").Value = 13848
$ws.Range("M132").Value = -7126.000100000001
$ws.Range("N132").Value = -18908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1666.6666
$ws.Range("I96").Value = 1500
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 1500
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -127
$ws.Range("N96").Value = -4746
$ws.Range("H132").Value = 1471.8148
$ws.Range("I132").Value = 928.94116
$ws.Range("J132").Value = 2394.7
$ws.Range("K132").Value = 2786.82348
$ws.Range("L132").Value = 7184.099999999999
$ws.Range("M132").Value = -256.82348
$ws.Range("N132").Value = -12244.1
$ws.Range("H136").Value = 23812168
$ws.Range("I136").Value = 30304210
$ws.Range("J136").Value = 8011.6665
$ws.Range("K136").Value = 90912630
$ws.Range("L136").Value = 24034.9995
$ws.Range("M136").Value = -90910080
$ws.Range("N136").Value = -29134.9995
